# Updated symbol list (crypto price/volume/hour refresh), matching commit
# "Updated symbol list on Sat Jan 14 13:07:28 UTC 2023 with GitHub Actions".
# Every data row (2-51) has its Hora (column G) bumped from 12 to 13, and
# most rows also get refreshed Price (D) and Volume(1h) (E) values.
# NumberFormat is forced to Text ("@") before each write so the values are
# stored as literal strings (matching the original inlineStr cells) instead
# of being auto-converted to numbers/percentages by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '304.05'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '6.41%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '13'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '8.64%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '13'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.269'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '3.75%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '13'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07560'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '12.65%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '13'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '7.15%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '13'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.756'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '9.03%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '13'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.491'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '8.82%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '13'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9145'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.11%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '13'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01665'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2,464.29%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '13'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1693'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '7.01%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '13'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07458'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '5.48%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '13'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08026'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '5.08%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '13'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.02997'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '2.49%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '13'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09899'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '10.21%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '13'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-6.18%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '13'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04551'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.98%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '13'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006283'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.44%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '13'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.495'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.31%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '13'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.232'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.05%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '13'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.42%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '13'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1341'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.05%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '13'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.492'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '12.91%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '13'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1626'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '4.15%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '13'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001213'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.82%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '13'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004446'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '1.79%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '13'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001326'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '10.28%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '13'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001736'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '7.21%'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '13'

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '13'

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '13'

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '13'

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '13'

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '13'

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '13'

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '13'

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '13'

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '13'

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '13'

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '13'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04514'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '6.50%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '13'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007208'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '6.51%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '13'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '8.99%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '13'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002245'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.49%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '13'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01304'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '2.81%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '13'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006202'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '7.49%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '13'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7091'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-63.88%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '13'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01296'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-13.65%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '13'

$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '13'

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '13'

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '13'

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '13'
